$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-04-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-22 Tuesday", 2) | Out-Null

# Update each arithmetic-problem cell in the table (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "12+54="
$t.Cell(1,2).Range.Text = "15+82="
$t.Cell(1,3).Range.Text = "24+45="
$t.Cell(1,4).Range.Text = "54+5="
$t.Cell(1,5).Range.Text = "16+20="

$t.Cell(2,1).Range.Text = "56+13="
$t.Cell(2,2).Range.Text = "44+8="
$t.Cell(2,3).Range.Text = "95-37="
$t.Cell(2,4).Range.Text = "65-6="
$t.Cell(2,5).Range.Text = "78-26="

$t.Cell(3,1).Range.Text = "5+28="
$t.Cell(3,2).Range.Text = "37+26="
$t.Cell(3,3).Range.Text = "52+45="
$t.Cell(3,4).Range.Text = "86-14="
$t.Cell(3,5).Range.Text = "96-32="

$t.Cell(4,1).Range.Text = "25+18="
$t.Cell(4,2).Range.Text = "3+64="
$t.Cell(4,3).Range.Text = "85-1="
$t.Cell(4,4).Range.Text = "6+62="
$t.Cell(4,5).Range.Text = "91-34="

$t.Cell(5,1).Range.Text = "1+70="
$t.Cell(5,2).Range.Text = "95-89="
$t.Cell(5,3).Range.Text = "47+21="
$t.Cell(5,4).Range.Text = "83-17="
$t.Cell(5,5).Range.Text = "77-27="

$t.Cell(6,1).Range.Text = "24-2="
$t.Cell(6,2).Range.Text = "32-20="
$t.Cell(6,3).Range.Text = "37+56="
$t.Cell(6,4).Range.Text = "38+37="
$t.Cell(6,5).Range.Text = "37-31="

$t.Cell(7,1).Range.Text = "68-36="
$t.Cell(7,2).Range.Text = "77-28="
$t.Cell(7,3).Range.Text = "73+20="
$t.Cell(7,4).Range.Text = "97-42="
$t.Cell(7,5).Range.Text = "33-5="

$t.Cell(8,1).Range.Text = "7+87="
$t.Cell(8,2).Range.Text = "39+42="
$t.Cell(8,3).Range.Text = "55-9="
$t.Cell(8,4).Range.Text = "49+29="
$t.Cell(8,5).Range.Text = "30-4="

$t.Cell(9,1).Range.Text = "48-6="
$t.Cell(9,2).Range.Text = "8+66="
$t.Cell(9,3).Range.Text = "6+61="
$t.Cell(9,4).Range.Text = "64-52="
$t.Cell(9,5).Range.Text = "60-47="

$t.Cell(10,1).Range.Text = "87-19="
$t.Cell(10,2).Range.Text = "80-64="
$t.Cell(10,3).Range.Text = "78-9="
$t.Cell(10,4).Range.Text = "93-56="
$t.Cell(10,5).Range.Text = "7+1="

$t.Cell(11,1).Range.Text = "88-35="
$t.Cell(11,2).Range.Text = "59-31="
$t.Cell(11,3).Range.Text = "98-16="
$t.Cell(11,4).Range.Text = "43+44="
$t.Cell(11,5).Range.Text = "88-5="

$t.Cell(12,1).Range.Text = "98-24="
$t.Cell(12,2).Range.Text = "63-60="
$t.Cell(12,3).Range.Text = "58+18="
$t.Cell(12,4).Range.Text = "69-6="
$t.Cell(12,5).Range.Text = "95-16="

$t.Cell(13,1).Range.Text = "49+4="
$t.Cell(13,2).Range.Text = "4+39="
$t.Cell(13,3).Range.Text = "81-17="
$t.Cell(13,4).Range.Text = "11+69="
$t.Cell(13,5).Range.Text = "70-41="

$t.Cell(14,1).Range.Text = "69-8="
$t.Cell(14,2).Range.Text = "12+72="
$t.Cell(14,3).Range.Text = "69-63="
$t.Cell(14,4).Range.Text = "71-54="
$t.Cell(14,5).Range.Text = "88+2="

$t.Cell(15,1).Range.Text = "42+37="
$t.Cell(15,2).Range.Text = "81+18="
$t.Cell(15,3).Range.Text = "11+64="
$t.Cell(15,4).Range.Text = "47+19="
$t.Cell(15,5).Range.Text = "49-25="

$t.Cell(16,1).Range.Text = "57-32="
$t.Cell(16,2).Range.Text = "33+20="
$t.Cell(16,3).Range.Text = "41+2="
$t.Cell(16,4).Range.Text = "45-20="
$t.Cell(16,5).Range.Text = "31+36="

$t.Cell(17,1).Range.Text = "56+32="
$t.Cell(17,2).Range.Text = "50+34="
$t.Cell(17,3).Range.Text = "87-28="
$t.Cell(17,4).Range.Text = "78-3="
$t.Cell(17,5).Range.Text = "11+35="

$t.Cell(18,1).Range.Text = "82+3="
$t.Cell(18,2).Range.Text = "52-32="
$t.Cell(18,3).Range.Text = "89-17="
$t.Cell(18,4).Range.Text = "62-18="
$t.Cell(18,5).Range.Text = "23+37="

$t.Cell(19,1).Range.Text = "29+7="
$t.Cell(19,2).Range.Text = "72-1="
$t.Cell(19,3).Range.Text = "5+9="
$t.Cell(19,4).Range.Text = "91-42="
$t.Cell(19,5).Range.Text = "74+20="

$t.Cell(20,1).Range.Text = "89-42="
$t.Cell(20,2).Range.Text = "11+69="
$t.Cell(20,3).Range.Text = "75-45="
$t.Cell(20,4).Range.Text = "8+39="
$t.Cell(20,5).Range.Text = "48+23="
